$wb = $excel.ActiveWorkbook

# --- Sheet: Vessel density (adipose) ---
# Edit vessel density SE for (Lijnen, 2006): 30 -> 37
$wsAdiposeDensity = $wb.Worksheets.Item("Vessel density (adipose)")
$wsAdiposeDensity.Range("C2").Value = 37
$wsAdiposeDensity.Range("C2").Select() | Out-Null

# --- Sheet: Vessel density (tumor) ---
# Add two new rows: Goel et al., 2013 (Nude mice) and Goel et al., 2013 (C57BL6/J mice)
$wsTumorDensity = $wb.Worksheets.Item("Vessel density (tumor)")

$wsTumorDensity.Range("A9").Value = "Goel et al., 2013 (Nude mice)"
$wsTumorDensity.Range("B9").Value = 292.45
$wsTumorDensity.Range("C9").Value = 28.64

$wsTumorDensity.Range("A10").Value = "Goel et al., 2013 (C57BL6/J mice)"
$wsTumorDensity.Range("B10").Value = 211.93
$wsTumorDensity.Range("C10").Value = 25.6

# Resize the table to include the new rows
$loTumorDensity = $wsTumorDensity.ListObjects.Item("Table4")
$loTumorDensity.Resize($wsTumorDensity.Range("A1:C10"))

# --- Sheet: Vessel size (tumor) selection state ---
$wsTumorSize = $wb.Worksheets.Item("Vessel size (tumor)")
$wsTumorSize.Range("A4:A5").Select() | Out-Null

# Activate Vessel density (tumor) as the final active/selected tab
# (this also clears the tabSelected flag previously on CBM (muscle))
$wsTumorDensity.Activate() | Out-Null
$wsTumorDensity.Range("B11").Select() | Out-Null

$wb.Save()
